$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'97.295.39"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +2.30%  "
$ws.Range('E2').Style = 'Normal'

# Row 3
$ws.Range('D3').Value = "'3.585.24"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +0.86%  "
$ws.Range('E3').Style = 'Normal'

# Row 4
$ws.Range('E4').Value = "'  +0.03%  "
$ws.Range('E4').Style = 'Normal'

# Row 5
$ws.Range('D5').Value = "'241.46"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +2.32%  "
$ws.Range('E5').Style = 'Normal'

# Row 6
$ws.Range('D6').Value = "'656.06"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +1.07%  "
$ws.Range('E6').Style = 'Normal'

# Row 7
$ws.Range('E7').Value = "'  +15.08%  "
$ws.Range('E7').Style = 'Normal'

# Row 8
$ws.Range('D8').Value = "'0.417"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  +4.99%  "
$ws.Range('E8').Style = 'Normal'

# Row 9
$ws.Range('B9').Value = 'USDC'
$ws.Range('C9').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D9').Value = "'0.999"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -0.05%  "
$ws.Range('E9').Style = 'Normal'

# Row 10
$ws.Range('B10').Value = 'Cardano'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D10').Value = "'1.05"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +4.74%  "
$ws.Range('E10').Style = 'Normal'

# Row 11
$ws.Range('D11').Value = "'3.584.19"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +0.90%  "
$ws.Range('E11').Style = 'Normal'

# Row 12
$ws.Range('D12').Value = "'44.31"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  +4.88%  "
$ws.Range('E12').Style = 'Normal'

# Row 13
$ws.Range('D13').Value = "'0.204"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  +0.76%  "
$ws.Range('E13').Style = 'Normal'

# Row 14
$ws.Range('D14').Value = "'6.41"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -0.95%  "
$ws.Range('E14').Style = 'Normal'

# Row 15
$ws.Range('D15').Value = "'97.315.17"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +2.43%  "
$ws.Range('E15').Style = 'Normal'

# Row 16
$ws.Range('D16').Value = "'4.252.58"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +0.88%  "
$ws.Range('E16').Style = 'Normal'

# Row 17
$ws.Range('D17').Value = "'0.0000258"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +2.41%  "
$ws.Range('E17').Style = 'Normal'

# Row 18
$ws.Range('D18').Value = "'8.14"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +5.08%  "
$ws.Range('E18').Style = 'Normal'

# Row 19
$ws.Range('D19').Value = "'3.579.93"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +0.93%  "
$ws.Range('E19').Style = 'Normal'

# Row 20
$ws.Range('D20').Value = "'12.75"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  +1.45%  "
$ws.Range('E20').Style = 'Normal'

# Row 21
$ws.Range('D21').Value = "'17.97"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  +0.78%  "
$ws.Range('E21').Style = 'Normal'

# Row 22
$ws.Range('E22').Value = "'  +9.41%  "
$ws.Range('E22').Style = 'Normal'

# Row 23
$ws.Range('D23').Value = "'3.50"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +0.94%  "
$ws.Range('E23').Style = 'Normal'

# Row 24
$ws.Range('D24').Value = "'513.87"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +1.98%  "
$ws.Range('E24').Style = 'Normal'

# Row 25
$ws.Range('D25').Value = "'0.0000204"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  +5.33%  "
$ws.Range('E25').Style = 'Normal'

# Row 26
$ws.Range('D26').Value = "'6.82"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  +1.62%  "
$ws.Range('E26').Style = 'Normal'

# Row 27
$ws.Range('D27').Value = "'95.45"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +0.42%  "
$ws.Range('E27').Style = 'Normal'

# Row 28
$ws.Range('E28').Value = "'  +4.36%  "
$ws.Range('E28').Style = 'Normal'

# Row 29
$ws.Range('D29').Value = "'3.778.53"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  +0.89%  "
$ws.Range('E29').Style = 'Normal'

# Row 30
$ws.Range('D30').Value = "'0.159"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  +12.74%  "
$ws.Range('E30').Style = 'Normal'

# Row 31
$ws.Range('E31').Value = "'  -0.89%  "
$ws.Range('E31').Style = 'Normal'

# Row 32
$ws.Range('D32').Value = "'11.76"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +3.21%  "
$ws.Range('E32').Style = 'Normal'

# Row 33
$ws.Range('D33').Value = "'0.996"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -0.36%  "
$ws.Range('E33').Style = 'Normal'

# Row 34
$ws.Range('D34').Value = "'0.183"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  +3.54%  "
$ws.Range('E34').Style = 'Normal'

# Row 35
$ws.Range('D35').Value = "'1.00"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  +0.21%  "
$ws.Range('E35').Style = 'Normal'

# Row 36
$ws.Range('D36').Value = "'31.63"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -0.16%  "
$ws.Range('E36').Style = 'Normal'

# Row 37
$ws.Range('D37').Value = "'619.44"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  +8.96%  "
$ws.Range('E37').Style = 'Normal'

# Row 38
$ws.Range('D38').Value = "'8.74"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  +4.25%  "
$ws.Range('E38').Style = 'Normal'

# Row 39
$ws.Range('D39').Value = "'0.566"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  +1.82%  "
$ws.Range('E39').Style = 'Normal'

# Row 40
$ws.Range('E40').Value = "'  +2.39%  "
$ws.Range('E40').Style = 'Normal'

# Row 41
$ws.Range('D41').Value = "'1.95"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  +11.96%  "
$ws.Range('E41').Style = 'Normal'

# Row 42
$ws.Range('D42').Value = "'0.153"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +2.59%  "
$ws.Range('E42').Style = 'Normal'

# Row 44
$ws.Range('D44').Value = "'0.922"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +2.78%  "
$ws.Range('E44').Style = 'Normal'

# Row 45
$ws.Range('D45').Value = "'5.96"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +5.91%  "
$ws.Range('E45').Style = 'Normal'

# Row 46
$ws.Range('D46').Value = "'0.0438"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +6.86%  "
$ws.Range('E46').Style = 'Normal'

# Row 47
$ws.Range('D47').Value = "'2.30"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  +0.51%  "
$ws.Range('E47').Style = 'Normal'

# Row 48
$ws.Range('D48').Value = "'23.65"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  +1.16%  "
$ws.Range('E48').Style = 'Normal'

# Row 49
$ws.Range('D49').Value = "'33.11"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  +0.06%  "
$ws.Range('E49').Style = 'Normal'

# Row 50
$ws.Range('D50').Value = "'8.47"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +5.14%  "
$ws.Range('E50').Style = 'Normal'

# Row 51
$ws.Range('E51').Value = "'  -0.16%  "
$ws.Range('E51').Style = 'Normal'
